# Add a new row (row 4) to the error-code table with a new entry:
#   A4 = 3
#   B4 = "+的處裡出現問題"  (leading "+" so it must be entered as text via
#         a quote-prefix, matching the "quotePrefix" cell style seen for
#         this new row)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (borders/alignment) of the row above (row 3) onto
# the new row's B cell before filling in values, so the new style entry
# picks up the quote-prefix flag in one shot.
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "'+的處裡出現問題"

$ws.Range("D9").Select()
